## The commit swaps the colour scheme carried by the deck's (single) real
## theme part with the colour scheme that was previously only used by the
## Notes Master ("Office Theme" palette), i.e. the slide master's theme
## becomes the stock "Office" colours instead of the custom "Integral"
## colours it had before.
##
## PowerPoint's object model exposes the 12 theme colour slots through
## ThemeColorScheme (Colors(1..12) -> dk1, lt1, dk2, lt2, accent1-6, hlink,
## folHlink), backed by an RGBColor whose .RGB is a BGR-packed long (the
## classic OLE RGB() layout: R | G<<8 | B<<16). Driving all twelve slots
## through that API rewrites the <a:clrScheme> of the presentation's
## theme in place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the "Office Theme" colours (dk1..folHlink), packed as
# BGR longs for the RGB property (0x00BBGGRR).
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
